# Weekly update for "Fruta, Vega Modelo de Temuco - Damasco":
# Two brand-new records are published at the top of the data block
# (rows 19-20), which pushes every existing record below them down by
# two rows (old row 21 -> new row 23, ... old row 75 -> new row 77).
# The two rows vacated by the shift (new rows 21-22) are re-filled with
# the data that used to live in rows 19-20 before this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

function Set-DataRow {
    param($RowNum, $Values)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$RowNum"
        $ws.Range($addr).Value = $Values[$i]
    }
}

# Make room for the 2 new records: insert 2 blank rows right after the
# header/existing row 20, shifting the old rows 21..75 down to 23..77.
$ws.Rows("21:22").Insert()

# New record written into row 19 (replaces the previous row-19 record,
# which is relocated to row 21 below).
Set-DataRow 19 @(10,"Vega Modelo de Temuco","La Araucanía",44910,9,"Fruta",100103,"Frutos de hueso (carozo)",100103003,"Damasco","Dina","Especial",125,17000,17000,17000,"`$/bandeja 10 kilos","Región de O'Higgins",1700,10)

# New record written into row 20 (replaces the previous row-20 record,
# which is relocated to row 22 below).
Set-DataRow 20 @(10,"Vega Modelo de Temuco","La Araucanía",44910,9,"Fruta",100103,"Frutos de hueso (carozo)",100103003,"Damasco","Dina","Primera",110,22000,22000,22000,"`$/bandeja 18 kilos","Región de O'Higgins",1222,18)

# Former row 19 now lives at row 21.
Set-DataRow 21 @(10,"Vega Modelo de Temuco","La Araucanía",44557,9,"Fruta",100103,"Frutos de hueso (carozo)",100103003,"Damasco","Dina","Primera",95,7000,7000,7000,"`$/bandeja 6 kilos","Provincia de San Felipe de Aconcagua",1167,6)

# Former row 20 now lives at row 22.
Set-DataRow 22 @(10,"Vega Modelo de Temuco","La Araucanía",44557,9,"Fruta",100103,"Frutos de hueso (carozo)",100103003,"Damasco","Dina","Primera",35,20000,20000,20000,"`$/caja 18 kilos","Provincia de Quillota",1111,18)

Write-Host "Applied weekly Damasco update: inserted 2 rows, refreshed rows 19-22."
